$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new "07-dec" date column before EJ ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("EJ").Insert()
$ws1.Range("EJ1").Value = "07-dec"
$ws1.Range("EJ2:EJ25").Value = "-"

# --- Sheet "Gaz": append new row for 2025-12-05 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A170").NumberFormat = "@"
$ws2.Range("A170").Value = "2025-12-05"
$ws2.Range("A170").Style = "Normal"
$ws2.Range("B170").Value = 25.965

# --- Sheet "CO2": append new row for 2025-12-05 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A170").NumberFormat = "@"
$ws3.Range("A170").Value = "2025-12-05"
$ws3.Range("A170").Style = "Normal"
$ws3.Range("B170").Value = 81.78
